# Replace the field construct ( { m:'doc.html'.fromHTMLURI() } written as a
# real Word field with fldChar begin/end and instrText runs ) with plain
# literal text runs spelling out the same characters, keeping the
# _GoBack bookmark in place between "doc.html" and "'.fromHTMLURI()".

$d = $word.ActiveDocument

# Locate the field so we do not have to hard-code character offsets.
$field = $d.Fields.Item(1)

# The field's overall range covers the begin mark, the field code
# (instrText runs) and the end mark / separator.
$fieldRange = $field.Code
$start = $fieldRange.Start - 1   # include the leading fldChar "begin"
$end = $fieldRange.End + 1       # include the trailing fldChar "end"

$range = $d.Range($start, $end)

$apos = "'"

$bodyXml = '<w:r><w:t>{</w:t></w:r>' `
  + '<w:r><w:t>m</w:t></w:r>' `
  + '<w:r><w:t>:</w:t></w:r>' `
  + '<w:r><w:t>' + $apos + '</w:t></w:r>' `
  + '<w:r><w:t>doc.html</w:t></w:r>' `
  + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
  + '<w:bookmarkEnd w:id="0"/>' `
  + '<w:r><w:t>' + $apos + '.fromHTMLURI()</w:t></w:r>' `
  + '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
  + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
  + '<pkg:xmlData>' `
  + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
  + '<w:body><w:p>' + $bodyXml + '</w:p></w:body>' `
  + '</w:document>' `
  + '</pkg:xmlData></pkg:part></pkg:package>'

$range.InsertXML($packageXml)
